# The commit swaps the contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: the deck's main theme (currently the "Integral" /
# "Red Violet" palette, used by the slide master + the presentation as a
# whole) becomes the plain "Office Theme" palette that used to live in
# theme1.xml (which was only ever wired up to the notes master).
#
# The PowerPoint object model only exposes a single Theme for the whole
# deck (SlideMaster.Theme / NotesMaster.Theme / HandoutMaster.Theme all
# resolve to the same ThemeColorScheme, the one backing the slide
# master / presentation theme relationship), so the reachable half of the
# swap is recoloring that shared theme to the "Office" values. We also
# try to rename the theme / color scheme to match; some hosts persist
# that, others treat Name as read-only, so the renames are best effort
# and wrapped so they can't abort the script.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# Index -> (slot, target "Office Theme" RRGGBB) using the standard
# MsoThemeColorSchemeIndex ordering: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink.
$targets = @(
    @{ Index = 1;  Slot = "dk1";      Hex = "000000" },
    @{ Index = 2;  Slot = "lt1";      Hex = "FFFFFF" },
    @{ Index = 3;  Slot = "dk2";      Hex = "44546A" },
    @{ Index = 4;  Slot = "lt2";      Hex = "E7E6E6" },
    @{ Index = 5;  Slot = "accent1";  Hex = "5B9BD5" },
    @{ Index = 6;  Slot = "accent2";  Hex = "ED7D31" },
    @{ Index = 7;  Slot = "accent3";  Hex = "A5A5A5" },
    @{ Index = 8;  Slot = "accent4";  Hex = "FFC000" },
    @{ Index = 9;  Slot = "accent5";  Hex = "4472C4" },
    @{ Index = 10; Slot = "accent6";  Hex = "70AD47" },
    @{ Index = 11; Slot = "hlink";    Hex = "0563C1" },
    @{ Index = 12; Slot = "folHlink"; Hex = "954F72" }
)

foreach ($t in $targets) {
    $hex = $t.Hex
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $rgb = $r -bor ($g -shl 8) -bor ($b -shl 16)
    $colorScheme.Colors($t.Index).RGB = $rgb
}

# Best-effort renames (no-ops on hosts that treat these as read-only).
try { $theme.Name = "Office Theme" } catch { }
try { $colorScheme.Name = "Office" } catch { }
try { $master.Design.Name = "Office Theme" } catch { }
